# expenditure_of_time.xlsx - add a new TODO row (row 14) to the "time costs" sheet,
# widen column G to fit the new note, and move the window/selection down to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time costs")

# --- widen column G so the new note text has room (was auto "best fit") ---
$ws.Columns.Item(7).ColumnWidth = 49.3

# --- new row 14: same date/"x" pattern as the rows above it ---
$ws.Range("A14").Value = 41716
$ws.Range("C14").Value = "x"
$ws.Range("D14").Value = "x"
$ws.Range("E14").Value = "x"
$ws.Range("F14").Value = "x"

$g14 = @"
Graue Vierecke anstatt "Verdeckt"
Spielaufgabe disabled Spielen button bei gegner
"@
$ws.Range("G14").Value = $g14
$ws.Range("G14").WrapText = $true
$ws.Range("G14").HorizontalAlignment = -4131
$ws.Range("G14").VerticalAlignment = -4108

$h14 = @"
Buttons nur auslösen, wenn Anfang & ende des toches drauf sind!
RandomEnemy (serverseitig!)
Bei SpielEnde Benachrichtigung & Ändern des SpielenButtons
Logo anzeigen (inapp & icon!)
Frage Buttons schrift zu klein nach Auswertung!
Beendete Spiele (letzte 5) in Sync mit liefern & in Hauptmenü anzeigen.
Duellanfragen werden u.U mehrmals im Hauptmenü angezeigt! (popUp)
tastaturinput-enter --> Aktion auf screen! (login/suchen..)
weiterbutton durch swipe ersetzen
login führt manchmal nicht zum home screen
aktualisieren buttons in RÜ & home in navigationbar
"@
$ws.Range("H14").Value = $h14
$ws.Range("H14").WrapText = $true

$ws.Rows.Item(14).RowHeight = 255

# --- move the view: scroll the frozen pane down and select the next empty cell ---
$null = $ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$null = $ws.Range("B14").Select()

# --- nudge the workbook window geometry to match (best effort) ---
$excel.ActiveWindow.Top = 135
$excel.ActiveWindow.Height = 9690
